function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws 'D2' '26.632.57'
Set-TextValue $ws 'E2' '  -0.10%  '
Set-TextValue $ws 'D3' '1.645.56'
Set-TextValue $ws 'E3' '  +0.71%  '
Set-TextValue $ws 'D4' '1.01'
Set-TextValue $ws 'E4' '  +0.30%  '
Set-TextValue $ws 'D5' '216.09'
Set-TextValue $ws 'E5' '  +1.34%  '
Set-TextValue $ws 'D6' '0.504'
Set-TextValue $ws 'E6' '  +0.36%  '
Set-TextValue $ws 'E7' '  +0.15%  '
Set-TextValue $ws 'E8' '  -0.17%  '
Set-TextValue $ws 'E9' '  +0.65%  '
Set-TextValue $ws 'D10' '19.37'
Set-TextValue $ws 'E10' '  +0.75%  '
Set-TextValue $ws 'D11' '0.0845'
Set-TextValue $ws 'E11' '  -0.03%  '
Set-TextValue $ws 'D12' '1.875.25'
Set-TextValue $ws 'E12' '  +0.70%  '
Set-TextValue $ws 'B13' 'Polkadot'
Set-TextValue $ws 'C13' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws 'D13' '4.23'
Set-TextValue $ws 'E13' '  +3.19%  '
Set-TextValue $ws 'B14' 'WrappedEther'
Set-TextValue $ws 'C14' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws 'D14' '1.628.17'
Set-TextValue $ws 'E14' '  -1.24%  '
Set-TextValue $ws 'D15' '0.536'
Set-TextValue $ws 'E15' '  +1.77%  '
Set-TextValue $ws 'D16' '66.33'
Set-TextValue $ws 'E16' '  +4.48%  '
Set-TextValue $ws 'D17' '26.691.81'
Set-TextValue $ws 'E17' '  +0.18%  '
Set-TextValue $ws 'E18' '  +1.47%  '
Set-TextValue $ws 'D19' '219.73'
Set-TextValue $ws 'E19' '  +0.33%  '
Set-TextValue $ws 'E20' '  +0.17%  '
Set-TextValue $ws 'E21' '  +2.08%  '
Set-TextValue $ws 'E22' '  +1.86%  '
Set-TextValue $ws 'D23' '9.56'
Set-TextValue $ws 'E23' '  +1.13%  '
Set-TextValue $ws 'D24' '2.11'
Set-TextValue $ws 'E24' '  +9.19%  '
Set-TextValue $ws 'D25' '147.11'
Set-TextValue $ws 'E25' '  -0.98%  '
Set-TextValue $ws 'E26' '  +0.27%  '
Set-TextValue $ws 'E27' '  +0.12%  '
Set-TextValue $ws 'D28' '7.16'
Set-TextValue $ws 'E28' '  +3.38%  '
Set-TextValue $ws 'D29' '15.95'
Set-TextValue $ws 'E29' '  +2.96%  '
Set-TextValue $ws 'D30' '0.0517'
Set-TextValue $ws 'E30' '  +1.43%  '
Set-TextValue $ws 'E31' '  +0.59%  '
Set-TextValue $ws 'D32' '3.40'
Set-TextValue $ws 'E32' '  +2.73%  '
Set-TextValue $ws 'D33' '3.07'
Set-TextValue $ws 'E33' '  +2.13%  '
Set-TextValue $ws 'D34' '1.283.34'
Set-TextValue $ws 'E34' '  +6.13%  '
Set-TextValue $ws 'D35' '1.55'
Set-TextValue $ws 'E35' '  +2.41%  '
Set-TextValue $ws 'D36' '0.0185'
Set-TextValue $ws 'E36' '  +7.41%  '
Set-TextValue $ws 'E37' '  +0.30%  '
Set-TextValue $ws 'B38' 'ImmutableX'
Set-TextValue $ws 'C38' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws 'D38' '0.529'
Set-TextValue $ws 'E38' '  +4.83%  '
Set-TextValue $ws 'B39' 'ARBITRUM'
Set-TextValue $ws 'C39' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws 'D39' '0.830'
Set-TextValue $ws 'E39' '  +2.23%  '
Set-TextValue $ws 'E40' '  +0.20%  '
Set-TextValue $ws 'D41' '0.808'
Set-TextValue $ws 'E41' '  +1.89%  '
Set-TextValue $ws 'E42' '  -2.05%  '
Set-TextValue $ws 'D43' '5.46'
Set-TextValue $ws 'E43' '  +0.38%  '
Set-TextValue $ws 'D44' '1.786.99'
Set-TextValue $ws 'E44' '  +0.94%  '
Set-TextValue $ws 'D45' '93.57'
Set-TextValue $ws 'E45' '  +0.13%  '
Set-TextValue $ws 'D46' '60.11'
Set-TextValue $ws 'E46' '  +9.66%  '
Set-TextValue $ws 'E47' '  +3.98%  '
Set-TextValue $ws 'D48' '0.0517'
Set-TextValue $ws 'E48' '  +0.69%  '
Set-TextValue $ws 'D49' '7.83'
Set-TextValue $ws 'E49' '  +1.68%  '
Set-TextValue $ws 'D50' '0.0982'
Set-TextValue $ws 'E50' '  +3.85%  '
Set-TextValue $ws 'E51' '  -0.60%  '
